# Add the missing "Tapa—Tartu—Koidula" line as a new row at the bottom of
# the "Lines detail" sheet (row 54), carrying over the same formatting as
# the other data rows, then leave the selection on the next empty row
# below it (A55) as recorded in the saved workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lines detail")

$srcRow = 52
$newRow = 54

# Copy formatting (font/style) from an existing fully-populated data row so
# the new row matches the rest of the table, including column G which is
# blank on the immediately preceding row.
$ws.Range("A$srcRow`:H$srcRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item($srcRow).RowHeight

$ws.Cells.Item($newRow, 1).Value = 109
$ws.Cells.Item($newRow, 2).Value = "Tapa—Tartu—Koidula"
$ws.Cells.Item($newRow, 3).Value = 1876
$ws.Cells.Item($newRow, 4).Value = "1520 mm"
$ws.Cells.Item($newRow, 5).Value = "Operational"
$ws.Cells.Item($newRow, 6).Value = 2011
$ws.Cells.Item($newRow, 7).Value = "Southern section rebuilt"
$ws.Cells.Item($newRow, 8).Value = $false

# Move the active selection to A55, matching the workbook's saved view.
$ws.Range("A55").Select()
